$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for the B2:K16 block after the edit (ifoCAST full series evaluation).
# Each row's values have shifted one column to the left (dropping the former first
# value), rows 2-6 gained a brand-new trailing value in column K, and the trailing
# cells that no longer have data (including B16) are cleared.
$data = @(
    @(-0.1220406607065741, -0.3824596977052969, 0.4529952984531944, 0.3529159385505728, -0.0008705675291234075, 0.1473068229033219, 0.3047528119074769, -0.1959066446185939, 0.08445123591687528, -0.08752362192813637),
    @(-0.4000138190652613, 0.4354411770932299, 0.3353618171906084, -0.01842468888908786, 0.1297527015433575, 0.2871986905475125, -0.2134607659785583, 0.06689711455691082, -0.1050777432881008, 0.2007847964582307),
    @(0.5235484960442361, 0.4234691361416146, 0.06968263006191837, 0.2178600204943637, 0.3753060094985187, -0.1253534470275521, 0.155004433507917, -0.01697042433709459, 0.2888921154092369, -0.3234378372780493),
    @(0.6020437696293982, 0.248257263549702, 0.3964346539821473, 0.5538806429863024, 0.0532211864602315, 0.3335790669957007, 0.161604209150689, 0.4674667488970205, -0.1448632037902657, 0.5435081814088458),
    @(1.160284491302352, 1.308461881734797, 1.465907870738952, 0.9652484142128814, 1.245606294748351, 1.073631436903339, 1.379493976649671, 0.7671640239623843, 1.455535409161496, 1.178844253737389),
    @(0.3644876457690405, 0.5219336347731955, 0.02127417824712469, 0.3016320587825939, 0.1296572009375822, 0.4355197406839137, -0.1768102120033725, 0.511561173195739, 0.2348700177716323, $null),
    @(0.5259015322863045, 0.0252420757602338, 0.305599956295703, 0.1336250984506913, 0.4394876381970228, -0.1728423144902634, 0.5155290707088481, 0.2388379152847414, $null, $null),
    @(0.1608821659304752, 0.4412400464659443, 0.2692651886209327, 0.5751277283672642, -0.03720222432002201, 0.6511691608790895, 0.3744780054549828, $null, $null, $null),
    @(0.2004338646102796, 0.028459006765268, 0.3343215465115995, -0.2780084061756867, 0.4103629790234248, 0.1336718235993181, $null, $null, $null, $null),
    @(-0.01687220848682837, 0.2889903312595031, -0.3233396214277831, 0.3650317637713285, 0.08834060834722172, $null, $null, $null, $null, $null),
    @(0.2221289093234493, -0.3902010433638369, 0.2981703418352746, 0.02147918641116785, $null, $null, $null, $null, $null, $null),
    @(-0.4197872457205535, 0.268584139478558, -0.00810701594554874, $null, $null, $null, $null, $null, $null, $null),
    @(0.2504334827489171, -0.02625767267518964, $null, $null, $null, $null, $null, $null, $null, $null),
    @(-0.04428949692388896, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = 2 + $j   # column B = 2
        $cell = $ws.Cells.Item($rowIndex, $colIndex)
        if ($null -eq $rowValues[$j]) {
            $cell.ClearContents()
        } else {
            $cell.Value = $rowValues[$j]
        }
    }
}
